# Applies the "marksheet" recalculation edit:
#  - Row 10 (No.) / Row 11 (Marking) / Row 12 (Total) summary values are
#    recomputed (handles float/partial-credit marking input without
#    clobbering the rest of the sheet).
#  - The redundant third "Student Ans / Correct Ans" block (columns G:H)
#    is removed.
#  - The per-question "Correct Ans" column (D) is folded away for most
#    rows and instead surfaced directly in column A so the sheet matches
#    the simplified two-column (Student Ans / Correct Ans) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 10: No. summary (Right / Wrong / Not Attempt / Max) ----
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").HorizontalAlignment = -4108   # xlCenter (re-assert so the
                                                # style's own xf entry is reused)
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

# ---- Row 11: Marking scheme (now numeric, handles the float case) ----
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# ---- Row 12: Total ----
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "55/112"

# ---- Drop the third (G:H) Student Ans / Correct Ans block entirely ----
$ws.Range("G15:H21").Clear()

# ---- D16 now shows the correct answer directly (was blank) ----
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").HorizontalAlignment = -4108   # xlCenter (reuse the style's xf)
$ws.Range("D16").Value = "Option A"

# ---- Clear the old per-row D:E "Student Ans / Correct Ans" pair for
#      rows 19-40 (superseded by column A below) ----
$ws.Range("D19:E40").Clear()

# ---- Column A (rows 19-40) becomes the "Correct Ans" column ----
$correctAnswers = @{
    19 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = "Option A"
    26 = "Option C"
    28 = "Option D"
    29 = "Option D"
    31 = "Option D"
    32 = "Option C"
    36 = "Option D"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}

foreach ($row in $correctAnswers.Keys) {
    $cell = $ws.Range("A$row")
    $cell.Value = $correctAnswers[$row]
    if ($row -eq 36) {
        $cell.Style = "incorrectStyle"
    } else {
        $cell.Style = "correctStyle"
    }
    $cell.HorizontalAlignment = -4108   # xlCenter (reuse the style's xf)
}
